$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The changed shape is "Rechteck 35" (id=36, creationId {A6C43494-5A24-0C46-B7CE-38CD76B734D4}),
# the 33rd shape on slide 2. Only its vertical position (y offset) moved,
# from 3287705 EMU to 3322874 EMU; x/width/height are unchanged.
$sh = $s.Shapes.Item(33)

# Shape.Top/Left are expressed in points (1 pt = 12700 EMU) and stored as
# single-precision floats internally, so converting EMU -> points and back
# truncates; add half an EMU before dividing so it rounds to the exact value.
$targetEmu = 3322874
$sh.Top = ($targetEmu + 0.5) / 12700
